$d = $word.ActiveDocument

$pairs = @(
    @("2025-07-21 Monday", "2025-07-22 Tuesday"),
    @("77÷4=", "38÷2="),
    @("77÷7=", "46÷3="),
    @("60÷9=", "96÷4="),
    @("33÷3=", "22÷4="),
    @("50÷5=", "86÷5="),
    @("57÷6=", "24÷9="),
    @("13÷6=", "34÷2="),
    @("18÷6=", "73÷3="),
    @("16÷6=", "69÷7="),
    @("32÷2=", "82÷7="),
    @("57÷3=", "74÷3="),
    @("14÷2=", "28÷6="),
    @("98÷7=", "87÷8="),
    @("37÷4=", "10÷2="),
    @("65÷3=", "18÷6="),
    @("76÷9=", "36÷4="),
    @("61÷2=", "38÷2="),
    @("79÷6=", "14÷4="),
    @("39÷9=", "70÷3="),
    @("67÷6=", "35÷8="),
    @("44÷7=", "37÷6="),
    @("19÷9=", "38÷2="),
    @("67÷2=", "49÷5="),
    @("42÷9=", "51÷4="),
    @("22÷8=", "52÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
